# Refresh cached Universalis market-price columns (H:N) on the per-job
# "Leve Profit" tables. Values below are the latest fetched averages;
# only the market-data columns change, source leve data (A:G) is untouched.
$wb = $excel.ActiveWorkbook

# --- ALC sheet ---
$ws = $wb.Worksheets.Item("ALC")
# Row 12: Don't Be So Tallow / Beeswax
$ws.Range("H12").Value = 183
$ws.Range("J12").Value = 199.5
$ws.Range("L12").Value = 199.5
$ws.Range("N12").Value = -539.5

# Row 111: An Eye for Healing / Grade 1 Dexterity Alkahest
$ws.Range("H111").Value = 2727.2727
$ws.Range("I111").Value = 1730.7142
$ws.Range("K111").Value = 5192.142599999999
$ws.Range("M111").Value = -2125.142599999999

# Row 112: Making Ends Meet / Superior Spiritbond Potion
$ws.Range("H112").Value = 1831.6786
$ws.Range("I112").Value = 1354.6
$ws.Range("J112").Value = 1935.3914
$ws.Range("K112").Value = 4063.8
$ws.Range("L112").Value = 5806.174199999999
$ws.Range("M112").Value = -2955.8
$ws.Range("N112").Value = -8022.174199999999

# Row 116: Growing Up / Growth Formula Kappa
$ws.Range("H116").Value = 7000
$ws.Range("I116").Value = 6000
$ws.Range("J116").Value = 9000
$ws.Range("K116").Value = 6000
$ws.Range("L116").Value = 9000
$ws.Range("M116").Value = -2558
$ws.Range("N116").Value = -15884

# Row 137: Cutting Edge of Culinary Quality / Magnesia Whetstone
$ws.Range("H137").Value = 3632.2273
$ws.Range("I137").Value = 3722.0527
$ws.Range("K137").Value = 11166.1581
$ws.Range("M137").Value = -8616.158100000001

# Row 138: All-night Crafting / Cunning Craftsman's Tisane
$ws.Range("H138").Value = 3812.9048
$ws.Range("I138").Value = 1655.1786
$ws.Range("J138").Value = 4891.768
$ws.Range("K138").Value = 4965.5358
$ws.Range("L138").Value = 14675.304
$ws.Range("M138").Value = 174.4642000000003
$ws.Range("N138").Value = -24955.304


# --- ARM sheet ---
$ws = $wb.Worksheets.Item("ARM")
# Row 4: Eyes Bigger than the Plate / Bronze Plate
$ws.Range("H4").Value = 199
$ws.Range("I4").Value = 199
$ws.Range("K4").Value = 199
$ws.Range("M4").Value = -83

# Row 5: The Alloyed Truth / Bronze Rivets
$ws.Range("H5").Value = 114.5
$ws.Range("I5").Value = 29.5
$ws.Range("K5").Value = 29.5
$ws.Range("M5").Value = 82.5

# Row 32: Ingot We Trust / Steel Ingot
$ws.Range("H32").Value = 15877168
$ws.Range("I32").Value = 15877168
$ws.Range("K32").Value = 15877168
$ws.Range("M32").Value = -15876881

# Row 63: Rivets Run through It / Mythrite Rivets
$ws.Range("H63").Value = 7932
$ws.Range("I63").Value = 4982.1665
$ws.Range("K63").Value = 4982.1665
$ws.Range("M63").Value = -4296.1665

# Row 66: A Riveting Revival (L) / Mythrite Rivets
$ws.Range("H66").Value = 7932
$ws.Range("I66").Value = 4982.1665
$ws.Range("K66").Value = 24910.8325
$ws.Range("M66").Value = -21478.8325

# Row 102: Smells of Rich Tama-hagane / Tama-hagane Ingot
$ws.Range("H102").Value = 1380
$ws.Range("I102").Value = 1380
$ws.Range("K102").Value = 1380
$ws.Range("M102").Value = 242

# Row 110: Scheduled Maintenance / Deepgold Ingot
$ws.Range("H110").Value = 7975.5
$ws.Range("J110").Value = 7975
$ws.Range("L110").Value = 7975
$ws.Range("N110").Value = -12065

# Row 122: Haste for High Durium / High Durium Nugget
$ws.Range("H122").Value = 3497.5
$ws.Range("I122").Value = 3497.5
$ws.Range("K122").Value = 10492.5
$ws.Range("M122").Value = -8042.5


# --- BSM sheet ---
$ws = $wb.Worksheets.Item("BSM")
# Row 4: Mending Fences / Bronze Rivets
$ws.Range("H4").Value = 114.5
$ws.Range("I4").Value = 29.5
$ws.Range("K4").Value = 29.5
$ws.Range("M4").Value = 85.5

# Row 20: Smelt and Dealt / Iron Ingot
$ws.Range("H20").Value = 1443.75
$ws.Range("I20").Value = 1450
$ws.Range("J20").Value = 1437.5
$ws.Range("K20").Value = 1450
$ws.Range("L20").Value = 1437.5
$ws.Range("M20").Value = -1203
$ws.Range("N20").Value = -1931.5

# Row 33: Mors Dagger / Steel Broadsword
$ws.Range("H33").Value = 2000
$ws.Range("I33").Value = 2000
$ws.Range("K33").Value = 2000
$ws.Range("M33").Value = -1664

# Row 94: High Steal / High Steel Nugget
$ws.Range("H94").Value = 2706.2
$ws.Range("I94").Value = 3051.4736
$ws.Range("J94").Value = 1612.8334
$ws.Range("K94").Value = 3051.4736
$ws.Range("L94").Value = 1612.8334
$ws.Range("M94").Value = -2600.4736
$ws.Range("N94").Value = -2514.8334

# Row 134: Ruthenium Supremium / Ruthenium Ingot
$ws.Range("H134").Value = 2556.0476
$ws.Range("I134").Value = 2272.4736
$ws.Range("K134").Value = 6817.4208
$ws.Range("M134").Value = -4282.4208


# --- CRP sheet ---
$ws = $wb.Worksheets.Item("CRP")
# Row 14: Citizens' Canes / Ash Radical
$ws.Range("H14").Value = 2402.75
$ws.Range("J14").Value = 2402.75
$ws.Range("L14").Value = 2402.75
$ws.Range("N14").Value = -2742.75

# Row 31: Wall Not Found / Walnut Lumber
$ws.Range("H31").Value = 1840.4242
$ws.Range("I31").Value = 1807.3462
$ws.Range("J31").Value = 1963.2858
$ws.Range("K31").Value = 1807.3462
$ws.Range("L31").Value = 1963.2858
$ws.Range("M31").Value = -1512.3462
$ws.Range("N31").Value = -2553.2858

# Row 34: Armoires of the Rich and Famous / Walnut Lumber
$ws.Range("H34").Value = 1840.4242
$ws.Range("I34").Value = 1807.3462
$ws.Range("J34").Value = 1963.2858
$ws.Range("K34").Value = 1807.3462
$ws.Range("L34").Value = 1963.2858
$ws.Range("M34").Value = -1605.3462
$ws.Range("N34").Value = -2367.2858

# Row 39: An Expected Tourney / Ash Cavalry Bow
$ws.Range("H39").Value = 14975
$ws.Range("I39").Value = 4760
$ws.Range("J39").Value = 32000
$ws.Range("K39").Value = 4760
$ws.Range("L39").Value = 32000
$ws.Range("M39").Value = -4369
$ws.Range("N39").Value = -32782

# Row 49: Bend It Like Durendaire / Ash Cavalry Bow
$ws.Range("H49").Value = 14975
$ws.Range("I49").Value = 4760
$ws.Range("J49").Value = 32000
$ws.Range("K49").Value = 4760
$ws.Range("L49").Value = 32000
$ws.Range("M49").Value = -4578
$ws.Range("N49").Value = -32364

# Row 58: You Do the Heavy Lifting / Mahogany Lumber
$ws.Range("H58").Value = 1626.3529
$ws.Range("I58").Value = 862.8148
$ws.Range("K58").Value = 862.8148
$ws.Range("M58").Value = -659.8148

# Row 100: Run Before They Walk / Pine Cane
$ws.Range("H100").Value = 100779.5
$ws.Range("J100").Value = 100779.5
$ws.Range("L100").Value = 100779.5
$ws.Range("N100").Value = -102943.5

# Row 136: Turali Quality / Dark Mahogany Lumber
$ws.Range("H136").Value = 1626.3529
$ws.Range("I136").Value = 862.8148
$ws.Range("K136").Value = 2588.4444
$ws.Range("M136").Value = -38.44439999999986


# --- GSM sheet ---
$ws = $wb.Worksheets.Item("GSM")
# Row 80: Needs More Prayerbell / Hardsilver Ingot
$ws.Range("H80").Value = 7586.724
$ws.Range("I80").Value = 6444.8945
$ws.Range("K80").Value = 6444.8945
$ws.Range("M80").Value = -5446.8945

# Row 83: With a Noise That Reaches Heaven (L) / Hardsilver Ingot
$ws.Range("H83").Value = 7586.724
$ws.Range("I83").Value = 6444.8945
$ws.Range("K83").Value = 32224.4725
$ws.Range("M83").Value = -27232.4725

# Row 109: You're My Wonderhall / Hematite Earrings of Healing
$ws.Range("H109").Value = 50000
$ws.Range("J109").Value = 50000
$ws.Range("L109").Value = 50000
$ws.Range("N109").Value = -52080

# Row 122: Awarding Academic Excellence / Ametrine
$ws.Range("H122").Value = 3607.182
$ws.Range("J122").Value = 4230.091
$ws.Range("L122").Value = 12690.273
$ws.Range("N122").Value = -17590.273

# Row 132: On Board for Lar / Lar Ingot
$ws.Range("H132").Value = 1721.0869
$ws.Range("I132").Value = 1454.75
$ws.Range("J132").Value = 3496.6667
$ws.Range("K132").Value = 4364.25
$ws.Range("L132").Value = 10490.0001
$ws.Range("M132").Value = -1834.25
$ws.Range("N132").Value = -15550.0001


# --- LTW sheet ---
$ws = $wb.Worksheets.Item("LTW")
# Row 20: Choke Hold / Hard Leather Choker
$ws.Range("H20").Value = 14285.714
$ws.Range("I20").Value = 10000
$ws.Range("J20").Value = 15000
$ws.Range("K20").Value = 10000
$ws.Range("L20").Value = 15000
$ws.Range("M20").Value = -9774
$ws.Range("N20").Value = -15452

# Row 61: Spelling Me Softly / Raptor Leather
$ws.Range("H61").Value = 1312.1482
$ws.Range("I61").Value = 976.087
$ws.Range("K61").Value = 976.087
$ws.Range("M61").Value = -774.087

# Row 113: Peace in Rest / Atrociraptor Leather
$ws.Range("H113").Value = 1312.1482
$ws.Range("I113").Value = 976.087
$ws.Range("K113").Value = 976.087
$ws.Range("M113").Value = 1193.913

# Row 122: Hell on Leather / Gaja Leather
$ws.Range("H122").Value = 6135.136
$ws.Range("I122").Value = 2813.6
$ws.Range("K122").Value = 8440.799999999999
$ws.Range("M122").Value = -5990.799999999999


# --- WVR sheet ---
$ws = $wb.Worksheets.Item("WVR")
# Row 52: Party Animals / Linen Deerstalker
$ws.Range("H52").Value = 22280
$ws.Range("I52").Value = 17850.25
$ws.Range("J52").Value = 39999
$ws.Range("K52").Value = 17850.25
$ws.Range("L52").Value = 39999
$ws.Range("M52").Value = -17624.25
$ws.Range("N52").Value = -40451

# Row 132: Comfy Cabins / Snow Cotton Cloth
$ws.Range("H132").Value = 1230.9642
$ws.Range("I132").Value = 1182.3334
$ws.Range("K132").Value = 3547.0002
$ws.Range("M132").Value = -1017.0002
